# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 41 (the sheet's single data
# table runs from row 2 to row 74), pushing the existing rows 41-74 down to
# 42-75 and extending the used range to A1:R75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 41; Excel shifts rows
# 41..74 down to 42..75 and the sheet dimension grows to A1:R75.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly record.
$ws.Cells.Item(41, 1).Value  = 4
$ws.Cells.Item(41, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(41, 3).Value  = "Los Lagos"
$ws.Cells.Item(41, 4).Value  = 45216
$ws.Cells.Item(41, 5).Value  = 10
$ws.Cells.Item(41, 6).Value  = 300000000
$ws.Cells.Item(41, 7).Value  = "Espárragos"
$ws.Cells.Item(41, 8).Value  = "Sin especificar"
$ws.Cells.Item(41, 9).Value  = "Primera"
$ws.Cells.Item(41, 10).Value = 400
$ws.Cells.Item(41, 11).Value = 1800
$ws.Cells.Item(41, 12).Value = 2000
$ws.Cells.Item(41, 13).Value = 1900
$ws.Cells.Item(41, 14).Value = "$/kilo"
$ws.Cells.Item(41, 15).Value = "Provincia de Linares"
$ws.Cells.Item(41, 16).Value = 1900
$ws.Cells.Item(41, 17).Value = 1
$ws.Cells.Item(41, 18).Value = "Hortaliza"
